$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.052.42'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.646.20'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.68'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5229'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2611'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06365'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('E10').Value = '  -1.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07665'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.65%  '
$ws.Range('D12').Value = '1.648.45'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.423'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').Value = '1.868.98'
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5549'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').Value = '0.0₅8267'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.98'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.20%  '
$ws.Range('D18').Value = '26.067.38'
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.726'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '188.55'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.243'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.18'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1220'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.419'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.386'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05946'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.267'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.400'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.397'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.661'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9980'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.393'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.755'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5619'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -6.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01611'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.850'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8565'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('D43').Value = '1.028.15'
$ws.Range('E43').Value = '  -8.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.24'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('D45').Value = '1.795.55'
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('E46').Value = '  +1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.80'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.085'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05151'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4219'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.51%  '
